$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 6874
$ws1.Range("F10").Value = 0
$ws1.Range("F11").Value = 11896
$ws1.Range("F12").Value = 12433
$ws1.Range("F13").Value = 1294
$ws1.Range("F14").Value = 0
$ws1.Range("F24").Value = 1475
$ws1.Range("F25").Value = 877
$ws1.Range("F26").Value = 10
$ws1.Range("F29").Value = 2926
$ws1.Range("F36").Value = 14
$ws1.Range("F40").Value = 263
$ws1.Range("F42").Value = 0
$ws1.Range("F44").Value = 142
$ws1.Range("F46").Value = 912

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6378

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 11896
$ws4.Range("F11").Value = 12433
$ws4.Range("F13").Value = 1294
$ws4.Range("F14").Value = 1264
$ws4.Range("F25").Value = 1475
$ws4.Range("F28").Value = 2926
$ws4.Range("F29").Value = 0
$ws4.Range("F34").Value = 6
$ws4.Range("F41").Value = 263
$ws4.Range("F44").Value = 142
$ws4.Range("F45").Value = 912
$ws4.Range("F46").Value = 274
